# Applies the commit "actualizacion a la lista de cosas por hacer":
#  1. Inserts a new "Sugerencia" paragraph (flanked by blank paragraphs)
#     right before the "REQUISITO 17 ... ESTADISTICAS" paragraph.
#  2. Moves the <w:lastRenderedPageBreak/> marker from the "La primer
#     categoría..." paragraph to the "un ID y una descripción..." paragraph.
#  3. Removes the stray <w:lastRenderedPageBreak/> marker from the
#     '///COPIA DE EL .TEXT "LO"' paragraph.
#
# NOTE: Range.InsertXML replaces the *entire* paragraph(s) touched by the
# range (even when the range is collapsed at a paragraph boundary), it
# does not splice new paragraphs in between neighbours. So every call
# below selects the full paragraph (Start..End, including the trailing
# paragraph mark) and re-supplies that paragraph's own content alongside
# any new paragraphs, instead of inserting at a collapsed point.

$d = $word.ActiveDocument
$wdPkg = 'xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"'

function New-Body-Xml([string]$bodyInner) {
    return '<pkg:package ' + $wdPkg + '><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Insert the new "Sugerencia" block before the REQUISITO 17 paragraph
#    (replace that paragraph with [blank, Sugerencia, blank, original]).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("REQUISITO 17 INCOMPLETO REVISAR EN DONDE PONER LAS ESTADISTICAS") | Out-Null
$target = $rng.Paragraphs(1).Range
$targetFull = $d.Range($target.Start, $target.End)

$sugerenciaBody = '<w:p/>' + `
    '<w:p>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Sugerencia</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">  :</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> ingresar estad&#237;sticas en &#8220;Admin ventas&#8221;</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p/>' + `
    '<w:p><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>REQUISITO 17 INCOMPLETO REVISAR EN DONDE PONER LAS ESTADISTICAS.</w:t></w:r></w:p>'

$targetFull.InsertXML((New-Body-Xml $sugerenciaBody))

# ---------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from "La primer categoría..." to
#    "un ID y una descripción... ese ID debe generarse".
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("un ID y una descripci") | Out-Null
$idPara = $rng2.Paragraphs(1).Range
$idParaFull = $d.Range($idPara.Start, $idPara.End)

$idBody = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>un ID y una descripci&#243;n&#8230; ese ID debe generarse</w:t></w:r></w:p>'
$idParaFull.InsertXML((New-Body-Xml $idBody))

$rng3 = $d.Content
$rng3.Find.Execute("primer categor") | Out-Null
$catPara = $rng3.Paragraphs(1).Range
$catParaFull = $d.Range($catPara.Start, $catPara.End)

$catBody = '<w:p>' + `
        '<w:r><w:t xml:space="preserve">La </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>primer categor&#237;a</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> que se ingresa en la BD ser&#225; la n&#250;mero</w:t></w:r>' + `
    '</w:p>'
$catParaFull.InsertXML((New-Body-Xml $catBody))

# ---------------------------------------------------------------------
# 3) Remove the stray <w:lastRenderedPageBreak/> before
#    '///COPIA DE EL .TEXT "LO"'.
# ---------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute("COPIA DE EL .TEXT") | Out-Null
$copiaPara = $rng4.Paragraphs(1).Range
$copiaParaFull = $d.Range($copiaPara.Start, $copiaPara.End)

$copiaBody = '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>///COPIA DE EL .TEXT &#8220;LO&#8221;</w:t></w:r>' + `
    '</w:p>'
$copiaParaFull.InsertXML((New-Body-Xml $copiaBody))

Write-Output "edits applied"
